# Add season-record columns (Wins / Losses / Ties) to the team sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/bordered header formatting used by the other
# header cells (copy format only from the neighboring header cell).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# The 2005 Yankees finished 95-67-0; stamp that record on every player row.
$wins = 95
$losses = 67
$ties = 0

for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
